$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 427, pushing the existing row 427
# (and everything below it) down by one row.
$ws.Rows(427).Insert()

# Populate the newly inserted row 427 with the new record's data.
$ws.Range("A427").Value = 4
$ws.Range("B427").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C427").Value = "Los Lagos"
$ws.Range("D427").Value = 45275
$ws.Range("E427").Value = 10
$ws.Range("F427").Value = 100112032
$ws.Range("G427").Value = "Zapallo italiano"
$ws.Range("H427").Value = "Sin especificar"
$ws.Range("I427").Value = "Primera"
$ws.Range("J427").Value = 250
$ws.Range("K427").Value = 16000
$ws.Range("L427").Value = 16000
$ws.Range("M427").Value = 16000
$ws.Range("N427").Value = "`$/caja 50 unidades"
$ws.Range("O427").Value = "Región de O'Higgins"
$ws.Range("P427").Value = 320
$ws.Range("Q427").Value = 50
$ws.Range("R427").Value = "Hortaliza"
